$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 94.5
$ws.Range("I12").Value = 94.5
$ws.Range("K12").Value = 94.5
$ws.Range("M12").Value = 75.5

$ws.Range("H19").Value = 701.6667
$ws.Range("J19").Value = 701.6667
$ws.Range("L19").Value = 701.6667
$ws.Range("N19").Value = -1051.6667

$ws.Range("H20").Value = 80000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H28").Value = 458.5
$ws.Range("I28").Value = 458.5
$ws.Range("K28").Value = 458.5
$ws.Range("M28").Value = 26.5

$ws.Range("H35").Value = 80000
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H41").Value = 675.44446
$ws.Range("I41").Value = 597.5
$ws.Range("J41").Value = 737.8
$ws.Range("K41").Value = 597.5
$ws.Range("L41").Value = 737.8
$ws.Range("M41").Value = -157.5
$ws.Range("N41").Value = -1617.8

$ws.Range("H98").Value = 2488.9167
$ws.Range("I98").Value = 857.9
$ws.Range("J98").Value = 10644
$ws.Range("K98").Value = 857.9
$ws.Range("L98").Value = 10644
$ws.Range("M98").Value = 640.1
$ws.Range("N98").Value = -13640

$ws.Range("H107").Value = 1635.0588
$ws.Range("I107").Value = 1612.25
$ws.Range("K107").Value = 1612.25
$ws.Range("M107").Value = 307.75

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H118").Value = 2284.6667
$ws.Range("J118").Value = 2656.9714
$ws.Range("L118").Value = 7970.914199999999
$ws.Range("N118").Value = -11284.9142

$ws.Range("H122").Value = 2488.9167
$ws.Range("I122").Value = 857.9
$ws.Range("J122").Value = 10644
$ws.Range("K122").Value = 2573.7
$ws.Range("L122").Value = 31932
$ws.Range("M122").Value = -123.6999999999998
$ws.Range("N122").Value = -36832

$ws.Range("H138").Value = 6500.7417
$ws.Range("J138").Value = 6638.724
$ws.Range("L138").Value = 19916.172
$ws.Range("N138").Value = -30196.172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1945.8572
$ws.Range("I2").Value = 1945.8572
$ws.Range("K2").Value = 1945.8572
$ws.Range("M2").Value = -1832.8572

$ws.Range("H4").Value = 95.8
$ws.Range("I4").Value = 69.75
$ws.Range("K4").Value = 69.75
$ws.Range("M4").Value = 46.25

$ws.Range("H32").Value = 8461.243
$ws.Range("I32").Value = 8461.243
$ws.Range("K32").Value = 8461.243
$ws.Range("M32").Value = -8174.243

$ws.Range("H116").Value = 1945.8572
$ws.Range("I116").Value = 1945.8572
$ws.Range("K116").Value = 1945.8572
$ws.Range("M116").Value = 348.1428000000001

$ws.Range("H122").Value = 4951.75
$ws.Range("I122").Value = 5597.6665
$ws.Range("K122").Value = 16792.9995
$ws.Range("M122").Value = -14342.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1945.8572
$ws.Range("I3").Value = 1945.8572
$ws.Range("K3").Value = 1945.8572
$ws.Range("M3").Value = -1831.8572

$ws.Range("H64").Value = 3722.6667
$ws.Range("J64").Value = 5482.6665
$ws.Range("L64").Value = 5482.6665
$ws.Range("N64").Value = -5932.6665

$ws.Range("H67").Value = 3722.6667
$ws.Range("J67").Value = 5482.6665
$ws.Range("L67").Value = 5482.6665
$ws.Range("N67").Value = -7042.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 861.6667
$ws.Range("I22").Value = 861.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 861.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -511.6667
$ws.Range("N22").ClearContents()

$ws.Range("H99").Value = 4300
$ws.Range("J99").Value = 4300
$ws.Range("L99").Value = 4300
$ws.Range("N99").Value = -7296

$ws.Range("H126").Value = 4300
$ws.Range("J126").Value = 4300
$ws.Range("L126").Value = 12900
$ws.Range("N126").Value = -17840

$ws.Range("H132").Value = 3254.875
$ws.Range("I132").Value = 1999.3334
$ws.Range("K132").Value = 5998.0002
$ws.Range("M132").Value = -3468.0002

$ws.Range("H134").Value = 2405.389
$ws.Range("I134").Value = 2053.7273
$ws.Range("J134").Value = 2958
$ws.Range("K134").Value = 6161.1819
$ws.Range("L134").Value = 8874
$ws.Range("M134").Value = -3626.1819
$ws.Range("N134").Value = -13944

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1594
$ws.Range("I5").Value = 1594
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4782
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4670
$ws.Range("N5").ClearContents()

$ws.Range("H33").Value = 115.57895
$ws.Range("I33").Value = 142
$ws.Range("J33").Value = 16.5
$ws.Range("K33").Value = 852
$ws.Range("L33").Value = 99
$ws.Range("M33").Value = -569
$ws.Range("N33").Value = -665

$ws.Range("H80").Value = 13942.857
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 13942.857
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H98").Value = 3035.2856
$ws.Range("I98").Value = 3053.4
$ws.Range("J98").Value = 2990
$ws.Range("K98").Value = 9160.200000000001
$ws.Range("L98").Value = 8970
$ws.Range("M98").Value = -7662.200000000001
$ws.Range("N98").Value = -11966

$ws.Range("H122").Value = 1845.6
$ws.Range("I122").Value = 996.3333
$ws.Range("J122").Value = 2209.5715
$ws.Range("K122").Value = 8966.9997
$ws.Range("L122").Value = 19886.1435
$ws.Range("M122").Value = -6516.9997
$ws.Range("N122").Value = -24786.1435

$ws.Range("H132").Value = 2646.4119
$ws.Range("I132").Value = 1299.8334
$ws.Range("J132").Value = 3380.9092
$ws.Range("K132").Value = 11698.5006
$ws.Range("L132").Value = 30428.1828
$ws.Range("M132").Value = -9168.500599999999
$ws.Range("N132").Value = -35488.1828

$ws.Range("H135").Value = 1594
$ws.Range("I135").Value = 1594
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 14346
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11811
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 2218.2856
$ws.Range("I136").Value = 2007
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 6021
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -921
$ws.Range("N136").Value = -17700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2647.6667
$ws.Range("I126").Value = 2647.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7943.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5473.000100000001
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3764.0715
$ws.Range("I132").Value = 2463.75
$ws.Range("K132").Value = 7391.25
$ws.Range("M132").Value = -4861.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376

$ws.Range("H136").Value = 4438.222
$ws.Range("I136").Value = 1927.8572
$ws.Range("K136").Value = 5783.571599999999
$ws.Range("M136").Value = -3233.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6526.8887
$ws.Range("I62").Value = 3033
$ws.Range("J62").Value = 8273.833000000001
$ws.Range("K62").Value = 3033
$ws.Range("L62").Value = 8273.833000000001
$ws.Range("M62").Value = -2409
$ws.Range("N62").Value = -9521.833000000001

$ws.Range("H65").Value = 6526.8887
$ws.Range("I65").Value = 3033
$ws.Range("J65").Value = 8273.833000000001
$ws.Range("K65").Value = 15165
$ws.Range("L65").Value = 41369.165
$ws.Range("M65").Value = -12045
$ws.Range("N65").Value = -47609.165

$ws.Range("H132").Value = 4402
$ws.Range("I132").Value = 1875.1428
$ws.Range("K132").Value = 5625.428400000001
$ws.Range("M132").Value = -3095.428400000001

$ws.Range("H136").Value = 983.63635
$ws.Range("I136").Value = 993.2
$ws.Range("K136").Value = 2979.6
$ws.Range("M136").Value = -429.6000000000004
